# Update the "cryptos" worksheet with refreshed price/volume(1h) figures,
# matching the GitHub Actions scheduled data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.470.82"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").Value = "2.031.93"
$ws.Range("E3").Value = "  +2.07%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'231.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -14.39%  "

$ws.Range("D6").Value = "'0.598"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.89%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'55.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("D10").Value = "'57.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("E11").Value = "  -2.30%  "

$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").Value = "2.328.67"
$ws.Range("E13").Value = "  +2.00%  "

$ws.Range("D14").Value = "'14.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").Value = "'20.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.36%  "

$ws.Range("D16").Value = "'0.759"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.58%  "

$ws.Range("D17").Value = "'5.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("D18").Value = "2.041.29"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("D19").Value = "36.783.78"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").Value = "'67.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.95%  "

$ws.Range("D21").Value = "'5.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.37%  "

$ws.Range("E22").Value = "  -4.07%  "

$ws.Range("D23").Value = "'220.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.80%  "

$ws.Range("D25").Value = "'2.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("D26").Value = "'2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.77%  "

$ws.Range("D27").Value = "'162.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("E28").Value = "  -2.69%  "

$ws.Range("D29").Value = "'0.127"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.84%  "

$ws.Range("D30").Value = "'18.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.95%  "

$ws.Range("D31").Value = "'1.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("E32").Value = "  -1.45%  "

$ws.Range("D33").Value = "'4.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.49%  "

# Rows 34/35: Hedera and LidoDAOToken swap rank positions
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'2.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.34%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.73%  "

$ws.Range("D36").Value = "'4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.19%  "

$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("E38").Value = "  -2.87%  "

$ws.Range("E39").Value = "  +8.76%  "

$ws.Range("E40").Value = "  -5.11%  "

$ws.Range("D41").Value = "'2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.57%  "

$ws.Range("D42").Value = "1.467.44"
$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").Value = "'0.0931"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.69%  "

# Rows 44/45: Aave and FTXToken swap rank positions
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +41.61%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'92.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.36%  "

$ws.Range("E46").Value = "  -2.04%  "

$ws.Range("E47").Value = "  -5.93%  "

$ws.Range("D48").Value = "'15.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "

$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").Value = "'6.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.27%  "

